$wb = $excel.ActiveWorkbook

# Sheet "requiremens" - add new requirement row
$ws1 = $wb.Worksheets.Item("requiremens")
$ws1.Range("A6").Value = "Implement ng-table for Definitions table"
$ws1.Range("A6").WrapText = $true

# Sheet "Bugs" - add Status value for third bug
$ws2 = $wb.Worksheets.Item("Bugs")
$ws2.Range("B3").Value = "Resolved"

# Update selections on both sheets to match the recorded cursor position
$ws2.Activate()
$ws2.Range("B4").Select()

$ws1.Activate()
$ws1.Range("A7").Select()
